$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 29, pushing the existing row 29 (and below) down.
$ws.Rows.Item(29).Insert()

# Populate the new row 29 with the new weekly entry.
$ws.Cells.Item(29, 1).Value = 10
$ws.Cells.Item(29, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(29, 3).Value = "La Araucanía"
$ws.Cells.Item(29, 4).Value = 45239
$ws.Cells.Item(29, 4).NumberFormat = $ws.Cells.Item(28, 4).NumberFormat
$ws.Cells.Item(29, 5).Value = 9
$ws.Cells.Item(29, 6).Value = "Fruta"
$ws.Cells.Item(29, 7).Value = 100104
$ws.Cells.Item(29, 8).Value = "Frutos de pepita"
$ws.Cells.Item(29, 9).Value = 100104004
$ws.Cells.Item(29, 10).Value = "Níspero"
$ws.Cells.Item(29, 11).Value = "Californiana(o)"
$ws.Cells.Item(29, 12).Value = "Primera"
$ws.Cells.Item(29, 13).Value = 55
$ws.Cells.Item(29, 14).Value = 26000
$ws.Cells.Item(29, 15).Value = 26000
$ws.Cells.Item(29, 16).Value = 26000
$ws.Cells.Item(29, 17).Value = '$/bandeja 5 kilos'
$ws.Cells.Item(29, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(29, 19).Value = 5200
$ws.Cells.Item(29, 20).Value = 5
